$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Inscritos (column E) 21 -> 22
$ws.Range("E10").Value = 22

# Row 16: Inscritos (E) 284 -> 286, Pagos (F) 79 -> 80, Inscrições homologadas (H) 79 -> 80
$ws.Range("E16").Value = 286
$ws.Range("F16").Value = 80
$ws.Range("H16").Value = 80
